$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Sheets.Item("ALC")
$ws.Cells.Item(93, 8).Value = 53248.5
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 53248.5
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 12).Value = 53248.5
$ws.Cells.Item(93, 14).Value = -58240.5
$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 14).Value = 0
$ws.Cells.Item(137, 12).ClearContents()
$ws.Cells.Item(137, 13).ClearContents()
$ws.Cells.Item(138, 8).Value = 5810.4443
$ws.Cells.Item(138, 9).Value = 1999
$ws.Cells.Item(138, 10).Value = 6286.875
$ws.Cells.Item(138, 11).Value = 5997
$ws.Cells.Item(138, 12).Value = 18860.625
$ws.Cells.Item(138, 13).Value = -857
$ws.Cells.Item(138, 14).Value = -29140.625

# ---- Sheet: ARM ----
$ws = $wb.Sheets.Item("ARM")
$ws.Cells.Item(19, 8).Value = 8399.6
$ws.Cells.Item(19, 9).Value = 8249.5
$ws.Cells.Item(19, 10).Value = 9000
$ws.Cells.Item(19, 11).Value = 8249.5
$ws.Cells.Item(19, 12).Value = 9000
$ws.Cells.Item(19, 13).Value = -8020.5
$ws.Cells.Item(19, 14).Value = -9458
$ws.Cells.Item(24, 8).Value = 2014141.8
$ws.Cells.Item(24, 9).Value = 0
$ws.Cells.Item(24, 10).Value = 2014141.8
$ws.Cells.Item(24, 11).Value = 0
$ws.Cells.Item(24, 12).Value = 2014141.8
$ws.Cells.Item(24, 14).Value = -2014889.8
$ws.Cells.Item(39, 8).Value = 6505.3335
$ws.Cells.Item(39, 9).Value = 4758
$ws.Cells.Item(39, 10).Value = 10000
$ws.Cells.Item(39, 11).Value = 4758
$ws.Cells.Item(39, 12).Value = 10000
$ws.Cells.Item(39, 13).Value = -4238
$ws.Cells.Item(39, 14).Value = -11040
$ws.Cells.Item(61, 8).Value = 3700.9443
$ws.Cells.Item(61, 9).Value = 3323.5334
$ws.Cells.Item(61, 10).Value = 5588
$ws.Cells.Item(61, 11).Value = 3323.5334
$ws.Cells.Item(61, 12).Value = 5588
$ws.Cells.Item(61, 13).Value = -3111.5334
$ws.Cells.Item(61, 14).Value = -6012
$ws.Cells.Item(74, 8).Value = 1318.3529
$ws.Cells.Item(74, 9).Value = 1318.3529
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 1318.3529
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 13).Value = -444.3529000000001
$ws.Cells.Item(77, 8).Value = 1318.3529
$ws.Cells.Item(77, 9).Value = 1318.3529
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 6591.7645
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 13).Value = -2223.7645
$ws.Cells.Item(100, 8).Value = 2014141.8
$ws.Cells.Item(100, 9).Value = 0
$ws.Cells.Item(100, 10).Value = 2014141.8
$ws.Cells.Item(100, 11).Value = 0
$ws.Cells.Item(100, 12).Value = 2014141.8
$ws.Cells.Item(100, 14).Value = -2016305.8
$ws.Cells.Item(112, 8).Value = 19685
$ws.Cells.Item(112, 9).Value = 0
$ws.Cells.Item(112, 10).Value = 19685
$ws.Cells.Item(112, 11).Value = 0
$ws.Cells.Item(112, 12).Value = 19685
$ws.Cells.Item(112, 14).Value = -22639
$ws.Cells.Item(132, 8).Value = 3968.6667
$ws.Cells.Item(132, 9).Value = 3968.6667
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 11906.0001
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 14).Value = -9376.000100000001
$ws.Cells.Item(132, 13).ClearContents()
$ws.Cells.Item(136, 8).Value = 3700.9443
$ws.Cells.Item(136, 9).Value = 3323.5334
$ws.Cells.Item(136, 10).Value = 5588
$ws.Cells.Item(136, 11).Value = 9970.600199999999
$ws.Cells.Item(136, 12).Value = 16764
$ws.Cells.Item(136, 13).Value = -7420.600199999999
$ws.Cells.Item(136, 14).Value = -21864

# ---- Sheet: BSM ----
$ws = $wb.Sheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 848.0476
$ws.Cells.Item(22, 9).Value = 754.2143
$ws.Cells.Item(22, 10).Value = 1035.7142
$ws.Cells.Item(22, 11).Value = 754.2143
$ws.Cells.Item(22, 12).Value = 1035.7142
$ws.Cells.Item(22, 13).Value = -581.2143
$ws.Cells.Item(22, 14).Value = -1381.7142
$ws.Cells.Item(76, 8).Value = 10000
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 10).Value = 10000
$ws.Cells.Item(76, 11).Value = 0
$ws.Cells.Item(76, 12).Value = 10000
$ws.Cells.Item(76, 14).Value = -10630
$ws.Cells.Item(79, 8).Value = 10000
$ws.Cells.Item(79, 9).Value = 0
$ws.Cells.Item(79, 10).Value = 10000
$ws.Cells.Item(79, 11).Value = 0
$ws.Cells.Item(79, 12).Value = 10000
$ws.Cells.Item(79, 14).Value = -12184
$ws.Cells.Item(111, 8).Value = 0
$ws.Cells.Item(111, 9).Value = 0
$ws.Cells.Item(111, 10).Value = 0
$ws.Cells.Item(111, 11).Value = 0
$ws.Cells.Item(111, 14).Value = 0
$ws.Cells.Item(111, 12).ClearContents()
$ws.Cells.Item(134, 8).Value = 3866.923
$ws.Cells.Item(134, 9).Value = 3843.6365
$ws.Cells.Item(134, 10).Value = 3995
$ws.Cells.Item(134, 11).Value = 11530.9095
$ws.Cells.Item(134, 12).Value = 11985
$ws.Cells.Item(134, 13).Value = -8995.9095
$ws.Cells.Item(134, 14).Value = -17055

# ---- Sheet: CRP ----
$ws = $wb.Sheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 883.3333
$ws.Cells.Item(16, 9).Value = 825
$ws.Cells.Item(16, 10).Value = 1000
$ws.Cells.Item(16, 11).Value = 825
$ws.Cells.Item(16, 12).Value = 1000
$ws.Cells.Item(16, 13).Value = -538
$ws.Cells.Item(16, 14).Value = -1574
$ws.Cells.Item(31, 8).Value = 3792.8386
$ws.Cells.Item(31, 9).Value = 1681.8948
$ws.Cells.Item(31, 10).Value = 4725.5815
$ws.Cells.Item(31, 11).Value = 1681.8948
$ws.Cells.Item(31, 12).Value = 4725.5815
$ws.Cells.Item(31, 13).Value = -1386.8948
$ws.Cells.Item(31, 14).Value = -5315.5815
$ws.Cells.Item(34, 8).Value = 3792.8386
$ws.Cells.Item(34, 9).Value = 1681.8948
$ws.Cells.Item(34, 10).Value = 4725.5815
$ws.Cells.Item(34, 11).Value = 1681.8948
$ws.Cells.Item(34, 12).Value = 4725.5815
$ws.Cells.Item(34, 13).Value = -1479.8948
$ws.Cells.Item(34, 14).Value = -5129.5815
$ws.Cells.Item(35, 8).Value = 1853.6666
$ws.Cells.Item(35, 9).Value = 1377.6666
$ws.Cells.Item(35, 10).Value = 3281.6667
$ws.Cells.Item(35, 11).Value = 1377.6666
$ws.Cells.Item(35, 12).Value = 3281.6667
$ws.Cells.Item(35, 13).Value = -1083.6666
$ws.Cells.Item(35, 14).Value = -3869.6667
$ws.Cells.Item(39, 8).Value = 0
$ws.Cells.Item(39, 9).Value = 0
$ws.Cells.Item(39, 10).Value = 0
$ws.Cells.Item(39, 11).Value = 0
$ws.Cells.Item(39, 12).Value = 0
$ws.Cells.Item(39, 13).ClearContents()
$ws.Cells.Item(41, 8).Value = 202410.62
$ws.Cells.Item(41, 9).Value = 6903.25
$ws.Cells.Item(41, 10).Value = 233691.8
$ws.Cells.Item(41, 11).Value = 6903.25
$ws.Cells.Item(41, 12).Value = 233691.8
$ws.Cells.Item(41, 13).Value = -6475.25
$ws.Cells.Item(41, 14).Value = -234547.8
$ws.Cells.Item(49, 8).Value = 0
$ws.Cells.Item(49, 9).Value = 0
$ws.Cells.Item(49, 10).Value = 0
$ws.Cells.Item(49, 11).Value = 0
$ws.Cells.Item(49, 12).Value = 0
$ws.Cells.Item(49, 13).ClearContents()
$ws.Cells.Item(50, 8).Value = 99888
$ws.Cells.Item(50, 9).Value = 0
$ws.Cells.Item(50, 10).Value = 99888
$ws.Cells.Item(50, 11).Value = 0
$ws.Cells.Item(50, 12).Value = 99888
$ws.Cells.Item(50, 14).Value = -101138
$ws.Cells.Item(51, 8).Value = 65644
$ws.Cells.Item(51, 9).Value = 10000
$ws.Cells.Item(51, 10).Value = 76772.8
$ws.Cells.Item(51, 11).Value = 10000
$ws.Cells.Item(51, 12).Value = 76772.8
$ws.Cells.Item(51, 13).Value = -9264
$ws.Cells.Item(51, 14).Value = -78244.8
$ws.Cells.Item(58, 8).Value = 1713.4615
$ws.Cells.Item(58, 9).Value = 934.63635
$ws.Cells.Item(58, 10).Value = 5997
$ws.Cells.Item(58, 11).Value = 934.63635
$ws.Cells.Item(58, 12).Value = 5997
$ws.Cells.Item(58, 13).Value = -731.63635
$ws.Cells.Item(58, 14).Value = -6403
$ws.Cells.Item(60, 8).Value = 62935.668
$ws.Cells.Item(60, 9).Value = 23000
$ws.Cells.Item(60, 10).Value = 94884.2
$ws.Cells.Item(60, 11).Value = 23000
$ws.Cells.Item(60, 12).Value = 94884.2
$ws.Cells.Item(60, 13).Value = -22489
$ws.Cells.Item(60, 14).Value = -95906.2
$ws.Cells.Item(61, 8).Value = 65644
$ws.Cells.Item(61, 9).Value = 10000
$ws.Cells.Item(61, 10).Value = 76772.8
$ws.Cells.Item(61, 11).Value = 10000
$ws.Cells.Item(61, 12).Value = 76772.8
$ws.Cells.Item(61, 13).Value = -9652
$ws.Cells.Item(61, 14).Value = -77468.8
$ws.Cells.Item(88, 8).Value = 8666.666999999999
$ws.Cells.Item(88, 9).Value = 6000
$ws.Cells.Item(88, 10).Value = 10000
$ws.Cells.Item(88, 11).Value = 6000
$ws.Cells.Item(88, 12).Value = 10000
$ws.Cells.Item(88, 13).Value = -5594
$ws.Cells.Item(88, 14).Value = -10812
$ws.Cells.Item(91, 8).Value = 8666.666999999999
$ws.Cells.Item(91, 9).Value = 6000
$ws.Cells.Item(91, 10).Value = 10000
$ws.Cells.Item(91, 11).Value = 6000
$ws.Cells.Item(91, 12).Value = 10000
$ws.Cells.Item(91, 13).Value = -4596
$ws.Cells.Item(91, 14).Value = -12808
$ws.Cells.Item(99, 8).Value = 1587
$ws.Cells.Item(99, 9).Value = 1264.75
$ws.Cells.Item(99, 10).Value = 2016.6666
$ws.Cells.Item(99, 11).Value = 1264.75
$ws.Cells.Item(99, 12).Value = 2016.6666
$ws.Cells.Item(99, 13).Value = 233.25
$ws.Cells.Item(99, 14).Value = -5012.6666
$ws.Cells.Item(107, 8).Value = 497.54544
$ws.Cells.Item(107, 9).Value = 399.42856
$ws.Cells.Item(107, 10).Value = 669.25
$ws.Cells.Item(107, 11).Value = 399.42856
$ws.Cells.Item(107, 12).Value = 669.25
$ws.Cells.Item(107, 13).Value = 1520.57144
$ws.Cells.Item(107, 14).Value = -4509.25
$ws.Cells.Item(113, 8).Value = 883.3333
$ws.Cells.Item(113, 9).Value = 825
$ws.Cells.Item(113, 10).Value = 1000
$ws.Cells.Item(113, 11).Value = 825
$ws.Cells.Item(113, 12).Value = 1000
$ws.Cells.Item(113, 13).Value = 1345
$ws.Cells.Item(113, 14).Value = -5340
$ws.Cells.Item(126, 8).Value = 1587
$ws.Cells.Item(126, 9).Value = 1264.75
$ws.Cells.Item(126, 10).Value = 2016.6666
$ws.Cells.Item(126, 11).Value = 3794.25
$ws.Cells.Item(126, 12).Value = 6049.9998
$ws.Cells.Item(126, 13).Value = -1324.25
$ws.Cells.Item(126, 14).Value = -10989.9998
$ws.Cells.Item(132, 8).Value = 1479.8
$ws.Cells.Item(132, 9).Value = 1479.8
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 4439.4
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -1909.4
$ws.Cells.Item(136, 8).Value = 1713.4615
$ws.Cells.Item(136, 9).Value = 934.63635
$ws.Cells.Item(136, 10).Value = 5997
$ws.Cells.Item(136, 11).Value = 2803.90905
$ws.Cells.Item(136, 12).Value = 17991
$ws.Cells.Item(136, 13).Value = -253.9090500000002
$ws.Cells.Item(136, 14).Value = -23091

# ---- Sheet: CUL ----
$ws = $wb.Sheets.Item("CUL")
$ws.Cells.Item(37, 8).Value = 175354.8
$ws.Cells.Item(37, 9).Value = 0
$ws.Cells.Item(37, 10).Value = 175354.8
$ws.Cells.Item(37, 11).Value = 0
$ws.Cells.Item(37, 12).Value = 526064.3999999999
$ws.Cells.Item(37, 14).Value = -526288.3999999999
$ws.Cells.Item(55, 8).Value = 4473.7
$ws.Cells.Item(55, 9).Value = 551
$ws.Cells.Item(55, 10).Value = 7088.8335
$ws.Cells.Item(55, 11).Value = 1653
$ws.Cells.Item(55, 12).Value = 21266.5005
$ws.Cells.Item(55, 13).Value = -1476
$ws.Cells.Item(55, 14).Value = -21620.5005

# ---- Sheet: GSM ----
$ws = $wb.Sheets.Item("GSM")
$ws.Cells.Item(12, 8).Value = 2000
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 2000
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 2000
$ws.Cells.Item(12, 14).Value = -2280
$ws.Cells.Item(63, 9).Value = 30000
$ws.Cells.Item(63, 10).Value = 30000
$ws.Cells.Item(63, 11).Value = 30000
$ws.Cells.Item(63, 12).Value = 30000
$ws.Cells.Item(63, 13).Value = -29314
$ws.Cells.Item(63, 14).Value = -31372
$ws.Cells.Item(66, 9).Value = 30000
$ws.Cells.Item(66, 10).Value = 30000
$ws.Cells.Item(66, 11).Value = 90000
$ws.Cells.Item(66, 12).Value = 90000
$ws.Cells.Item(66, 13).Value = -86568
$ws.Cells.Item(66, 14).Value = -96864
$ws.Cells.Item(92, 8).Value = 200251
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 200251
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 200251
$ws.Cells.Item(92, 14).Value = -203995
$ws.Cells.Item(132, 8).Value = 2857.0908
$ws.Cells.Item(132, 9).Value = 2857.0908
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 8571.2724
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -6041.2724

# ---- Sheet: LTW ----
$ws = $wb.Sheets.Item("LTW")
$ws.Cells.Item(33, 8).Value = 1625
$ws.Cells.Item(33, 9).Value = 1400
$ws.Cells.Item(33, 10).Value = 1850
$ws.Cells.Item(33, 11).Value = 1400
$ws.Cells.Item(33, 12).Value = 1850
$ws.Cells.Item(33, 13).Value = -1110
$ws.Cells.Item(33, 14).Value = -2430
$ws.Cells.Item(46, 8).Value = 4517.16
$ws.Cells.Item(46, 9).Value = 3190.7693
$ws.Cells.Item(46, 10).Value = 5954.0835
$ws.Cells.Item(46, 11).Value = 3190.7693
$ws.Cells.Item(46, 12).Value = 5954.0835
$ws.Cells.Item(46, 13).Value = -3002.7693
$ws.Cells.Item(46, 14).Value = -6330.0835
$ws.Cells.Item(127, 8).Value = 0
$ws.Cells.Item(127, 9).Value = 0
$ws.Cells.Item(127, 10).Value = 0
$ws.Cells.Item(127, 11).Value = 0
$ws.Cells.Item(127, 14).Value = 0
$ws.Cells.Item(127, 12).ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Sheets.Item("WVR")
$ws.Cells.Item(69, 8).Value = 38387.855
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 38387.855
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 12).Value = 38387.855
$ws.Cells.Item(69, 14).Value = -39885.855
$ws.Cells.Item(72, 8).Value = 38387.855
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 10).Value = 38387.855
$ws.Cells.Item(72, 11).Value = 0
$ws.Cells.Item(72, 12).Value = 115163.565
$ws.Cells.Item(72, 14).Value = -122651.565
$ws.Cells.Item(92, 8).Value = 25000
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 25000
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 25000
$ws.Cells.Item(92, 14).Value = -29992
$ws.Cells.Item(113, 8).Value = 475.1
$ws.Cells.Item(113, 9).Value = 358.83334
$ws.Cells.Item(113, 10).Value = 649.5
$ws.Cells.Item(113, 11).Value = 1076.50002
$ws.Cells.Item(113, 12).Value = 1948.5
$ws.Cells.Item(113, 13).Value = 1093.49998
$ws.Cells.Item(113, 14).Value = -6288.5
$ws.Cells.Item(126, 8).Value = 6199.6665
$ws.Cells.Item(126, 9).Value = 4624
$ws.Cells.Item(126, 10).Value = 6987.5
$ws.Cells.Item(126, 11).Value = 13872
$ws.Cells.Item(126, 12).Value = 20962.5
$ws.Cells.Item(126, 13).Value = -11402
$ws.Cells.Item(126, 14).Value = -25902.5
$ws.Cells.Item(132, 8).Value = 1764.9166
$ws.Cells.Item(132, 9).Value = 1598.1818
$ws.Cells.Item(132, 10).Value = 3599
$ws.Cells.Item(132, 11).Value = 4794.5454
$ws.Cells.Item(132, 12).Value = 10797
$ws.Cells.Item(132, 13).Value = -2264.5454
$ws.Cells.Item(132, 14).Value = -15857
